$wb = $excel.ActiveWorkbook

# --- Data sheet updates (ordered first so new shared strings line up) ---
$data = $wb.Worksheets.Item("Data")
# A3: "September STEO" -> "November STEO"
$data.Range("A3").Value = "November STEO"
# Updated GDP values
$data.Range("B3").Value = 19092
$data.Range("C3").Value = 18411
$data.Range("D3").Value = 19098

# --- About sheet updates ---
$about = $wb.Worksheets.Item("About")
# B6: "January 2020 and September 2020" -> "January 2020 and November 2020"
$about.Range("B6").Value = "January 2020 and November 2020"
# A27: "As of EPS 2.1.1, this variable is set up to model the impacts of the 2020" -> "As of EPS 3.1, ..."
$about.Range("A27").Value = "As of EPS 3.1, this variable is set up to model the impacts of the 2020"
# A28: "SARS-CoV-2 pandemic.  It uses the latest data available as of September 9," -> "...November 10,"
$about.Range("A28").Value = "SARS-CoV-2 pandemic.  It uses the latest data available as of November 10,"

# Update selection to match diff (D4 -> B12), while keeping "About" the active sheet
$data.Range("B12").Select()
$about.Activate()

$wb.Save()
